$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New row 21 (write new shared strings in the same order the diff introduces them)
$ws.Range("D21").Value = "10 - Listar usuarios"
$ws.Range("F21").Value = "Método (GET)"
$ws.Range("G21").Value = "app.get(urls.users,isUser, listUsers)"
$ws.Range("H21").Value = "isUser"
$ws.Range("I21").Value = "sólo me sale para un campo"
$ws.Range("J21").Value = "ok"

# New cell on existing row 20: I20 (reuses the shared string added above)
$ws.Range("I20").Value = "sólo me sale para un campo"

# Update view: topLeftCell and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("I22").Select()

# Tab ratio on workbook window
$excel.ActiveWindow.TabRatio = 384
